$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in results for existing rows 96 and 105 ---
$ws.Cells.Item(96, 7).Value = "Fallo"
$ws.Cells.Item(96, 8).Value = -1

$ws.Cells.Item(105, 7).Value = "Fallo"
$ws.Cells.Item(105, 8).Value = -1

# --- Append new rows 106-111 ---
$newRows = @(
    @{ Row = 106; A = 14762056; B = "2025-10-07"; C = "Valentin Vacherot"; D = "Tallon Griekspoor"; E = "Gana Tallon Griekspoor"; F = 1.62 },
    @{ Row = 107; A = 14833288; B = "2025-10-06"; C = "Cannon Kingsley";   D = "Andrew Fenty";       E = "Gana Cannon Kingsley";   F = 2 },
    @{ Row = 108; A = 14831275; B = "2025-10-06"; C = "Andre Ilagan";     D = "Kaylan Bigun";        E = "Gana Kaylan Bigun";      F = 2.75 },
    @{ Row = 109; A = 14831044; B = "2025-10-06"; C = "Matias Soto";      D = "Miguel Tobon";        E = "Gana Miguel Tobon";      F = 2 },
    @{ Row = 110; A = 14832957; B = "2025-10-06"; C = "Ivan Gakhov";      D = "Martin Krumich";      E = "Gana Martin Krumich";    F = 1.73 },
    @{ Row = 111; A = 14832956; B = "2025-10-06"; C = "Sumit Nagal";      D = "Alexander Ritschard"; E = "Gana Alexander Ritschard"; F = 1.91 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A

    # Force column B to be stored as text so date-like strings are not
    # auto-converted into date serial values.
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $r.B

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F

    # resultado / profit remain blank (pending) for these newly appended
    # matches, matching the source data - intentionally left unset.
}
